$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "560.73") need their
# number format forced to Text first, otherwise Excel auto-converts the
# typed string into a numeric value - these price cells must stay text.
$textForceCells = @("D5", "D6", "D11", "D13", "D14", "D19", "D20", "D21", "D22", "D24", "D28", "D29", "D32", "D33", "D34", "D36", "D37", "D38", "D41", "D42", "D44", "D45", "D46", "D47", "D49")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.019.06"
$ws.Range("E2").Value = "  +2.76%  "
$ws.Range("D3").Value = "2.413.09"
$ws.Range("E3").Value = "  +3.99%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "560.73"
$ws.Range("E5").Value = "  +2.82%  "
$ws.Range("D6").Value = "138.53"
$ws.Range("E6").Value = "  +5.74%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +1.34%  "
$ws.Range("D9").Value = "2.413.40"
$ws.Range("E9").Value = "  +4.05%  "
$ws.Range("E10").Value = "  +3.44%  "
$ws.Range("D11").Value = "5.73"
$ws.Range("E11").Value = "  +4.30%  "
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "0.348"
$ws.Range("E13").Value = "  +3.85%  "
$ws.Range("D14").Value = "25.69"
$ws.Range("E14").Value = "  +8.55%  "
$ws.Range("D15").Value = "2.847.16"
$ws.Range("E15").Value = "  +4.19%  "
$ws.Range("D16").Value = "62.023.75"
$ws.Range("E16").Value = "  +2.93%  "
$ws.Range("E17").Value = "  +5.08%  "
$ws.Range("D18").Value = "2.420.62"
$ws.Range("E18").Value = "  +4.43%  "
$ws.Range("D19").Value = "11.09"
$ws.Range("E19").Value = "  +4.71%  "
$ws.Range("D20").Value = "344.54"
$ws.Range("E20").Value = "  +9.83%  "
$ws.Range("D21").Value = "4.22"
$ws.Range("E21").Value = "  +2.11%  "
$ws.Range("D22").Value = "6.86"
$ws.Range("E22").Value = "  +3.62%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").Value = "65.05"
$ws.Range("E24").Value = "  +1.98%  "
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("E27").Value = "  +6.94%  "
$ws.Range("D28").Value = "1.52"
$ws.Range("E28").Value = "  +13.24%  "
$ws.Range("D29").Value = "1.37"
$ws.Range("E29").Value = "  +15.79%  "
$ws.Range("E30").Value = "  +3.55%  "
$ws.Range("D31").Value = "0.0₃0781"
$ws.Range("E31").Value = "  +7.11%  "
$ws.Range("D32").Value = "6.36"
$ws.Range("E32").Value = "  +7.39%  "
$ws.Range("D33").Value = "170.45"
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("D34").Value = "1.41"
$ws.Range("E34").Value = "  +2.53%  "
$ws.Range("E35").Value = "  +3.62%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").Value = "376.26"
$ws.Range("E36").Value = "  +16.17%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "18.58"
$ws.Range("E37").Value = "  +4.18%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "4.53"
$ws.Range("E38").Value = "  +12.03%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "1.66"
$ws.Range("E41").Value = "  +9.19%  "
$ws.Range("D42").Value = "39.15"
$ws.Range("E42").Value = "  +3.20%  "
$ws.Range("E43").Value = "  +4.50%  "
$ws.Range("D44").Value = "3.66"
$ws.Range("E44").Value = "  +4.92%  "
$ws.Range("D45").Value = "20.50"
$ws.Range("E45").Value = "  +6.63%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0530"
$ws.Range("E46").Value = "  +6.94%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.0960"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("E48").Value = "  +4.82%  "
$ws.Range("D49").Value = "18.01"
$ws.Range("E49").Value = "  +6.97%  "
$ws.Range("E50").Value = "  +3.68%  "
$ws.Range("D51").Value = "0.0₆0214"
$ws.Range("E51").Value = "  +2.11%  "
